# Reproduce the target edit:
#  - Keep the physically-original worksheet (it owns the existing drawing
#    relationship) but rename it to "Sheet2" and wipe its data.
#  - Add a brand-new worksheet, rename it "Sheet1", and give it the new data
#    (shifted into column M, except the "nombre2" label which stays in C5),
#    reusing the existing cell style (so no duplicate style/font entries are
#    created) and the 5-5.8ish column width band used for the helper columns.

$wb = $excel.ActiveWorkbook
$orig = $wb.Worksheets.Item(1)

# Grab a styled source cell before we touch anything, so we can stamp the
# same style onto the new cells later without minting new style/font ids.
$styleSource = $orig.Range("A1")

# Create the replacement sheet right after the original one so the tab
# order ends up [orig, newSheet].
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $orig)

# Avoid a transient name collision ("Sheet2" is already the auto-assigned
# name of $newSheet at this point) by routing through a temp name.
$newSheet.Name = "__tmp_new__"
$orig.Name = "Sheet2"
$newSheet.Name = "Sheet1"

# --- Wipe the old data off the (now) "Sheet2" tab -------------------------
$orig.Range("A1").Clear()
$orig.Range("C2:C5").Clear()

# --- Populate the (now) "Sheet1" tab with the new data --------------------
$newSheet.Range("M2").Value = "a,1"
$newSheet.Range("M3").Value = "b,1"
$newSheet.Range("M4").Value = "s,3"
$newSheet.Range("C5").Value = "nombre2"
$newSheet.Range("M5").Value = "d,2"

# Re-apply the shared cell style (font/alignment) used by the original data
# cells so the new cells carry the same style index instead of a fresh one.
$styleSource.Copy()
$newSheet.Range("M2:M5").PasteSpecial(-4122)
$newSheet.Range("C5").PasteSpecial(-4122)

# --- Column widths: columns B:M (2-13) get a narrow custom width ----------
$newSheet.Columns("B:M").ColumnWidth = 4.166666666666667
